$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 45456070
$ws.Range("I125").Value = 90910504
$ws.Range("J125").Value = 1635.6364
$ws.Range("K125").Value = 818194536
$ws.Range("L125").Value = 14720.7276
$ws.Range("M125").Value = -818192076
$ws.Range("N125").Value = -19640.7276

$ws.Range("H126").Value = 30000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 30000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -39880

$ws.Range("H127").Value = 26453.896
$ws.Range("I127").Value = 77258.234
$ws.Range("J127").Value = 1051.7307
$ws.Range("K127").Value = 231774.702
$ws.Range("L127").Value = 3155.1921
$ws.Range("M127").Value = -226814.702
$ws.Range("N127").Value = -13075.1921

$ws.Range("H128").Value = 33000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 33000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 33000
$ws.Range("N128").Value = -42960

$ws.Range("H129").Value = 701.0741
$ws.Range("I129").Value = 529.3889
$ws.Range("J129").Value = 1044.4445
$ws.Range("K129").Value = 1588.1667
$ws.Range("L129").Value = 3133.3335
$ws.Range("M129").Value = 3411.8333
$ws.Range("N129").Value = -13133.3335

$ws.Range("H130").Value = 29600
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 29600
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 29600
$ws.Range("N130").Value = -39640

$ws.Range("H131").Value = 7431.5264
$ws.Range("I131").Value = 516.5833
$ws.Range("J131").Value = 19285.715
$ws.Range("K131").Value = 1549.7499
$ws.Range("L131").Value = 57857.145
$ws.Range("M131").Value = 3490.2501
$ws.Range("N131").Value = -67937.145

$ws.Range("H132").Value = 12170480
$ws.Range("I132").Value = 1354.9683
$ws.Range("J132").Value = 50503224
$ws.Range("K132").Value = 4064.9049
$ws.Range("L132").Value = 151509672
$ws.Range("M132").Value = -1534.9049
$ws.Range("N132").Value = -151514732

$ws.Range("H133").Value = 60000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 60000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120

$ws.Range("H134").Value = 60000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 60000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -70140

$ws.Range("H135").Value = 11112215
$ws.Range("I135").Value = 289.68
$ws.Range("J135").Value = 25002122
$ws.Range("K135").Value = 2607.12
$ws.Range("L135").Value = 225019098
$ws.Range("M135").Value = -72.11999999999989
$ws.Range("N135").Value = -225024168

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0

$ws.Range("H137").Value = 21768.125
$ws.Range("I137").Value = 30291.559
$ws.Range("J137").Value = 1068.3572
$ws.Range("K137").Value = 90874.677
$ws.Range("L137").Value = 3205.0716
$ws.Range("M137").Value = -88324.677
$ws.Range("N137").Value = -8305.0716

$ws.Range("H138").Value = 1598.71
$ws.Range("I138").Value = 902.0323
$ws.Range("J138").Value = 1911.7102
$ws.Range("K138").Value = 2706.0969
$ws.Range("L138").Value = 5735.1306
$ws.Range("M138").Value = 2433.9031
$ws.Range("N138").Value = -16015.1306

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H140").Value = 60000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 60000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360

$ws.Range("H141").Value = 2237.2222
$ws.Range("I141").Value = 1361.25
$ws.Range("J141").Value = 3332.1875
$ws.Range("K141").Value = 4083.75
$ws.Range("L141").Value = 9996.5625
$ws.Range("M141").Value = 1096.25
$ws.Range("N141").Value = -20356.5625


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 996317.25
$ws.Range("I105").Value = 1593207.8
$ws.Range("J105").Value = 1499.6666
$ws.Range("K105").Value = 1593207.8
$ws.Range("L105").Value = 1499.6666
$ws.Range("M105").Value = -1591460.8
$ws.Range("N105").Value = -4993.6666

$ws.Range("H133").Value = 56560
$ws.Range("J133").Value = 56560
$ws.Range("L133").Value = 56560
$ws.Range("N133").Value = -66680


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1507.5714
$ws.Range("I132").Value = 1217.25
$ws.Range("J132").Value = 3249.5
$ws.Range("K132").Value = 3651.75
$ws.Range("L132").Value = 9748.5
$ws.Range("M132").Value = -1121.75
$ws.Range("N132").Value = -14808.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 44823.22
$ws.Range("J129").Value = 68449.6
$ws.Range("L129").Value = 205348.8
$ws.Range("N129").Value = -215348.8

$ws.Range("H131").Value = 14205383
$ws.Range("J131").Value = 15625870
$ws.Range("L131").Value = 46877610
$ws.Range("N131").Value = -46887690

$ws.Range("H137").Value = 8972518
$ws.Range("I137").Value = 27778966
$ws.Range("J137").Value = 4458970
$ws.Range("K137").Value = 83336898
$ws.Range("L137").Value = 13376910
$ws.Range("M137").Value = -83331798
$ws.Range("N137").Value = -13387110


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 37608.25
$ws.Range("I52").Value = 19000
$ws.Range("J52").Value = 56216.5
$ws.Range("K52").Value = 19000
$ws.Range("L52").Value = 56216.5
$ws.Range("M52").Value = -18741
$ws.Range("N52").Value = -56734.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0

$ws.Range("H121").Value = 32000
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 32000
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 32000
$ws.Range("N121").Value = -35494

$ws.Range("H122").Value = 8600.77
$ws.Range("I122").Value = 4980
$ws.Range("J122").Value = 10863.75
$ws.Range("K122").Value = 14940
$ws.Range("L122").Value = 32591.25
$ws.Range("M122").Value = -12490
$ws.Range("N122").Value = -37491.25

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0

$ws.Range("H126").Value = 860.13635
$ws.Range("I126").Value = 601.2778
$ws.Range("J126").Value = 2025
$ws.Range("K126").Value = 1803.8334
$ws.Range("L126").Value = 6075
$ws.Range("M126").Value = 666.1666
$ws.Range("N126").Value = -11015

$ws.Range("H127").Value = 50000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 50000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws.Range("H129").Value = 48000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 48000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 48000
$ws.Range("N129").Value = -58000

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0

$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0

$ws.Range("H132").Value = 2749.2156
$ws.Range("I132").Value = 533.09753
$ws.Range("J132").Value = 11835.3
$ws.Range("K132").Value = 1599.29259
$ws.Range("L132").Value = 35505.89999999999
$ws.Range("M132").Value = 930.70741
$ws.Range("N132").Value = -40565.89999999999

$ws.Range("H133").Value = 40692.145
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 40692.145
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 40692.145
$ws.Range("N133").Value = -50812.145

$ws.Range("H135").Value = 38667.75
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 38667.75
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 38667.75
$ws.Range("N135").Value = -48807.75

$ws.Range("H136").Value = 1510383.5
$ws.Range("I136").Value = 1459769.5
$ws.Range("J136").Value = 2006401
$ws.Range("K136").Value = 4379308.5
$ws.Range("L136").Value = 6019203
$ws.Range("M136").Value = -4376758.5
$ws.Range("N136").Value = -6024303

$ws.Range("H137").Value = 53286.11
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 53286.11
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 53286.11
$ws.Range("N137").Value = -63486.11

$ws.Range("H138").Value = 50975
$ws.Range("I138").Value = 50000
$ws.Range("J138").Value = 51300
$ws.Range("K138").Value = 50000
$ws.Range("L138").Value = 51300
$ws.Range("M138").Value = -44860
$ws.Range("N138").Value = -61580

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H140").Value = 42500
$ws.Range("I140").Value = 20000
$ws.Range("J140").Value = 50000
$ws.Range("K140").Value = 20000
$ws.Range("L140").Value = 50000
$ws.Range("M140").Value = -14820
$ws.Range("N140").Value = -60360

$ws.Range("H141").Value = 54750
$ws.Range("I141").Value = 40000
$ws.Range("J141").Value = 56857.145
$ws.Range("K141").Value = 40000
$ws.Range("L141").Value = 56857.145
$ws.Range("M141").Value = -34820
$ws.Range("N141").Value = -67217.14499999999

